$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# rf (risk-free rate) adjusted to match the same return period used by the
# other series, instead of always using the annual rf as before. This
# changes GRS, p-value of GRS, A|a|, A|a|/A|re|, and A(a^2)/A(re^2) for
# part1's size-inv portfolios.
$pvalue = 0.000000000000000111022302462515703346169743122319869483035330

$ws.Cells.Item(2, 2).Value = 12.18812099985124
$ws.Cells.Item(2, 3).Value = $pvalue
$ws.Cells.Item(2, 4).Value = 0.00736770806979629
$ws.Cells.Item(2, 5).Value = 0.8929020919339536
$ws.Cells.Item(2, 6).Value = 0.7972741457800305

$ws.Cells.Item(3, 2).Value = 12.52358522724342
$ws.Cells.Item(3, 3).Value = $pvalue
$ws.Cells.Item(3, 4).Value = 0.007355251178049704
$ws.Cells.Item(3, 5).Value = 0.8913924250749858
$ws.Cells.Item(3, 6).Value = 0.7945804554810642

$ws.Cells.Item(4, 2).Value = 13.61217423570494
$ws.Cells.Item(4, 3).Value = $pvalue
$ws.Cells.Item(4, 4).Value = 0.006486985216250005
$ws.Cells.Item(4, 5).Value = 0.78616614760829
$ws.Cells.Item(4, 6).Value = 0.6180572116452596

$ws.Cells.Item(5, 2).Value = 13.8912692044324
$ws.Cells.Item(5, 3).Value = $pvalue
$ws.Cells.Item(5, 4).Value = 0.006812886420499159
$ws.Cells.Item(5, 5).Value = 0.825662537025618
$ws.Cells.Item(5, 6).Value = 0.68171862504758

$ws.Cells.Item(6, 2).Value = 13.81768963282079
$ws.Cells.Item(6, 3).Value = $pvalue
$ws.Cells.Item(6, 4).Value = 0.006478520402170357
$ws.Cells.Item(6, 5).Value = 0.7851402858168145
$ws.Cells.Item(6, 6).Value = 0.6164452684125092
